$d = $word.ActiveDocument

# --- Insert a new "Database" sub-bullet (ListParagraph / ilvl=2 /
# numId=2, same as the other Functional-requirement bullets) right
# after "Performance/Response time" and right before the "App"
# paragraph. ---
$find = $d.Content
[void]$find.Find.Execute("Performance/Response time")
$find.InsertParagraphAfter()

$newPara = $find.Paragraphs(1).Next()
$newPara.Range.Text = "DatabaseX"

# --- Move the hidden "_GoBack" bookmark from the end of the
# "Fearghal - 2 minutes" paragraph to the end of the new "Database"
# paragraph. ---
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

# Bookmarking a genuinely zero-length range is unreliable here, so
# bookmark the trailing placeholder character ("X") and then delete
# it - the bookmark collapses to sit right after "Database", matching
# the target layout.
$dbRange = $d.Content
[void]$dbRange.Find.Execute("DatabaseX")
$lastChar = $dbRange.Duplicate
$lastChar.Start = $lastChar.End - 1
$d.Bookmarks.Add("_GoBack", $lastChar)
$lastChar.Text = ""
